# Generate Report for Handback
# Inserts a new handback record for 55b57630-96f8-4b35-86f4-4d813023a6b7.md
# (status "Handed back: in sync with en-US") ahead of the existing
# 8d64a134-a970-43ce-bd7d-af5ac30db086.md record on all three sheets:
#   Sheet1 "Overview" (A1:G)
#   Sheet2 "zh-cn"     (A1:P)
#   Sheet3 "de-de"     (A1:P)
#
# Net effect (matches the OOXML diff): a new row 3 is written with the
# 55b57630 data, and the data that used to live in row 3 (8d64a134) is
# pushed down into a brand-new row 4.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).Value = $text
}

function Set-DateCell($ws, $addr, $text) {
    $ws.Range($addr).NumberFormat = $dateFmt
    $ws.Range($addr).Value = $text
}

function Set-LinkCell($ws, $addr, $text) {
    $ws.Range($addr).Style = "Hyperlink"
    $ws.Range($addr).Value = $text
}

# ---------------------------------------------------------------------
# Sheet1 : Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$lo1 = $ws1.ListObjects.Item(1)
$lo1.ListRows.Add() | Out-Null

# Row 3 -> new 55b57630 record
Set-TextCell $ws1 "A3" "55b57630-96f8-4b35-86f4-4d813023a6b7.md"
Set-LinkCell $ws1 "B3" "e2e\55b57630-96f8-4b35-86f4-4d813023a6b7.md"
Set-TextCell $ws1 "C3" ".md"
Set-TextCell $ws1 "E3" "Handed back: in sync with en-US"
Set-TextCell $ws1 "F3" "Handed back: in sync with en-US"
Set-DateCell $ws1 "G3" "2016-08-23 10:46:26"

# Row 4 -> the 8d64a134 record that used to sit in row 3
Set-TextCell $ws1 "A4" "8d64a134-a970-43ce-bd7d-af5ac30db086.md"
Set-LinkCell $ws1 "B4" "e2e\8d64a134-a970-43ce-bd7d-af5ac30db086.md"
Set-TextCell $ws1 "C4" ".md"
Set-TextCell $ws1 "E4" "Handed back: in sync with en-US"
Set-TextCell $ws1 "F4" "Handed back: in sync with en-US"
Set-DateCell $ws1 "G4" "2016-08-23 10:41:22"

# Rebuild the hyperlinks (Range.Hyperlinks.Delete() removes every
# hyperlink on the sheet, so recreate all of them in final order).
$ws1.Range("B2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edc3ce23a54c34bbaf4eb31160f6fd2fa98d5f78/e2e/8a45673f-20d6-4ebb-bb2b-8fae5966994c.md", "", "", "e2e\8a45673f-20d6-4ebb-bb2b-8fae5966994c.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/eeedb5cf7094a0684f30fbc1085e7372bf0d40a3/e2e/55b57630-96f8-4b35-86f4-4d813023a6b7.md", "", "", "e2e\55b57630-96f8-4b35-86f4-4d813023a6b7.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ca72dfe75ee4c7cf625fb7c17929cb45e2b02fe0/e2e/8d64a134-a970-43ce-bd7d-af5ac30db086.md", "", "", "e2e\8d64a134-a970-43ce-bd7d-af5ac30db086.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet2 : zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$lo2 = $ws2.ListObjects.Item(1)
$lo2.ListRows.Add() | Out-Null

# Row 3 -> new 55b57630 record
Set-LinkCell $ws2 "A3" "55b57630-96f8-4b35-86f4-4d813023a6b7.md"
Set-TextCell $ws2 "B3" ".md"
Set-TextCell $ws2 "C3" "Handed back: in sync with en-US"
Set-TextCell $ws2 "D3" "e2e"
Set-TextCell $ws2 "E3" "ht"
Set-TextCell $ws2 "F3" "True"
Set-TextCell $ws2 "G3" "55b57630-96f8-4b35-86f4-4d813023a6b7.eeedb5cf7094a0684f30fbc1085e7372bf0d40a3.zh-cn.xlf"
Set-DateCell $ws2 "H3" "2016-08-23 10:46:22"
Set-LinkCell $ws2 "I3" "55b57630-96f8-4b35-86f4-4d813023a6b7.md"
Set-TextCell $ws2 "J3" "55b57630-96f8-4b35-86f4-4d813023a6b7.eeedb5cf7094a0684f30fbc1085e7372bf0d40a3.zh-cn.xlf"
Set-DateCell $ws2 "K3" "2016-08-23 10:46:47"
Set-TextCell $ws2 "L3" ""
Set-TextCell $ws2 "M3" "True"
Set-TextCell $ws2 "N3" ""
Set-TextCell $ws2 "O3" "False"
Set-TextCell $ws2 "P3" ""

# Row 4 -> the 8d64a134 record that used to sit in row 3
Set-LinkCell $ws2 "A4" "8d64a134-a970-43ce-bd7d-af5ac30db086.md"
Set-TextCell $ws2 "B4" ".md"
Set-TextCell $ws2 "C4" "Handed back: in sync with en-US"
Set-TextCell $ws2 "D4" "e2e"
Set-TextCell $ws2 "E4" "ht"
Set-TextCell $ws2 "F4" "True"
Set-TextCell $ws2 "G4" "8d64a134-a970-43ce-bd7d-af5ac30db086.f8f53b79bc1fa8daaded2b2686b6dfbe2d2a839b.zh-cn.xlf"
Set-DateCell $ws2 "H4" "2016-08-23 10:41:17"
Set-LinkCell $ws2 "I4" "8d64a134-a970-43ce-bd7d-af5ac30db086.md"
Set-TextCell $ws2 "J4" "8d64a134-a970-43ce-bd7d-af5ac30db086.f8f53b79bc1fa8daaded2b2686b6dfbe2d2a839b.zh-cn.xlf"
Set-DateCell $ws2 "K4" "2016-08-23 10:41:35"
Set-TextCell $ws2 "L4" ""
Set-TextCell $ws2 "M4" "True"
Set-TextCell $ws2 "N4" ""
Set-TextCell $ws2 "O4" "False"
Set-TextCell $ws2 "P4" ""

$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edc3ce23a54c34bbaf4eb31160f6fd2fa98d5f78/e2e/8a45673f-20d6-4ebb-bb2b-8fae5966994c.md", "", "", "8a45673f-20d6-4ebb-bb2b-8fae5966994c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/29aad4fae15b6961ff83b9c77b9d976e83029ce1/e2e/8a45673f-20d6-4ebb-bb2b-8fae5966994c.md", "", "", "8a45673f-20d6-4ebb-bb2b-8fae5966994c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edc3ce23a54c34bbaf4eb31160f6fd2fa98d5f78/e2e/55b57630-96f8-4b35-86f4-4d813023a6b7.md", "", "", "55b57630-96f8-4b35-86f4-4d813023a6b7.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/eeedb5cf7094a0684f30fbc1085e7372bf0d40a3/e2e/55b57630-96f8-4b35-86f4-4d813023a6b7.md", "", "", "55b57630-96f8-4b35-86f4-4d813023a6b7.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ca72dfe75ee4c7cf625fb7c17929cb45e2b02fe0/e2e/8d64a134-a970-43ce-bd7d-af5ac30db086.md", "", "", "8d64a134-a970-43ce-bd7d-af5ac30db086.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fd4bbaacfc9b5b25c09f2a6cc8f4f83e2a10693f/e2e/8d64a134-a970-43ce-bd7d-af5ac30db086.md", "", "", "8d64a134-a970-43ce-bd7d-af5ac30db086.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet3 : de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$lo3 = $ws3.ListObjects.Item(1)
$lo3.ListRows.Add() | Out-Null

# Row 3 -> new 55b57630 record
Set-LinkCell $ws3 "A3" "55b57630-96f8-4b35-86f4-4d813023a6b7.md"
Set-TextCell $ws3 "B3" ".md"
Set-TextCell $ws3 "C3" "Handed back: in sync with en-US"
Set-TextCell $ws3 "D3" "e2e"
Set-TextCell $ws3 "E3" "ht"
Set-TextCell $ws3 "F3" "True"
Set-TextCell $ws3 "G3" "55b57630-96f8-4b35-86f4-4d813023a6b7.eeedb5cf7094a0684f30fbc1085e7372bf0d40a3.de-de.xlf"
Set-DateCell $ws3 "H3" "2016-08-23 10:46:26"
Set-LinkCell $ws3 "I3" "55b57630-96f8-4b35-86f4-4d813023a6b7.md"
Set-TextCell $ws3 "J3" "55b57630-96f8-4b35-86f4-4d813023a6b7.eeedb5cf7094a0684f30fbc1085e7372bf0d40a3.de-de.xlf"
Set-DateCell $ws3 "K3" "2016-08-23 10:46:54"
Set-TextCell $ws3 "L3" ""
Set-TextCell $ws3 "M3" "True"
Set-TextCell $ws3 "N3" ""
Set-TextCell $ws3 "O3" "False"
Set-TextCell $ws3 "P3" ""

# Row 4 -> the 8d64a134 record that used to sit in row 3
Set-LinkCell $ws3 "A4" "8d64a134-a970-43ce-bd7d-af5ac30db086.md"
Set-TextCell $ws3 "B4" ".md"
Set-TextCell $ws3 "C4" "Handed back: in sync with en-US"
Set-TextCell $ws3 "D4" "e2e"
Set-TextCell $ws3 "E4" "ht"
Set-TextCell $ws3 "F4" "True"
Set-TextCell $ws3 "G4" "8d64a134-a970-43ce-bd7d-af5ac30db086.f8f53b79bc1fa8daaded2b2686b6dfbe2d2a839b.de-de.xlf"
Set-DateCell $ws3 "H4" "2016-08-23 10:41:22"
Set-LinkCell $ws3 "I4" "8d64a134-a970-43ce-bd7d-af5ac30db086.md"
Set-TextCell $ws3 "J4" "8d64a134-a970-43ce-bd7d-af5ac30db086.f8f53b79bc1fa8daaded2b2686b6dfbe2d2a839b.de-de.xlf"
Set-DateCell $ws3 "K4" "2016-08-23 10:41:42"
Set-TextCell $ws3 "L4" ""
Set-TextCell $ws3 "M4" "True"
Set-TextCell $ws3 "N4" ""
Set-TextCell $ws3 "O4" "False"
Set-TextCell $ws3 "P4" ""

$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edc3ce23a54c34bbaf4eb31160f6fd2fa98d5f78/e2e/8a45673f-20d6-4ebb-bb2b-8fae5966994c.md", "", "", "8a45673f-20d6-4ebb-bb2b-8fae5966994c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0e78413e414bfa920e2cbf12663a97efc0fb176c/e2e/8a45673f-20d6-4ebb-bb2b-8fae5966994c.md", "", "", "8a45673f-20d6-4ebb-bb2b-8fae5966994c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edc3ce23a54c34bbaf4eb31160f6fd2fa98d5f78/e2e/55b57630-96f8-4b35-86f4-4d813023a6b7.md", "", "", "55b57630-96f8-4b35-86f4-4d813023a6b7.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/eeedb5cf7094a0684f30fbc1085e7372bf0d40a3/e2e/55b57630-96f8-4b35-86f4-4d813023a6b7.md", "", "", "55b57630-96f8-4b35-86f4-4d813023a6b7.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ca72dfe75ee4c7cf625fb7c17929cb45e2b02fe0/e2e/8d64a134-a970-43ce-bd7d-af5ac30db086.md", "", "", "8d64a134-a970-43ce-bd7d-af5ac30db086.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/aaf7dc677c7d51a07ae3728a9585f7701046f1e3/e2e/8d64a134-a970-43ce-bd7d-af5ac30db086.md", "", "", "8d64a134-a970-43ce-bd7d-af5ac30db086.md") | Out-Null
